$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "65.328.29"
$c.Style = $origStyle
$ws.Range("E2").Value = "  +3.68%  "
$c = $ws.Range("D3")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "3.487.14"
$c.Style = $origStyle
$ws.Range("E3").Value = "  +3.08%  "
$ws.Range("E4").Value = "  +0.03%  "
$c = $ws.Range("D5")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "579.40"
$c.Style = $origStyle
$ws.Range("E5").Value = "  +2.71%  "
$c = $ws.Range("D6")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "162.34"
$c.Style = $origStyle
$ws.Range("E6").Value = "  +4.74%  "
$ws.Range("E7").Value = "  +13.43%  "
$c = $ws.Range("D8")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = $origStyle
$ws.Range("E8").Value = "  +0.04%  "
$c = $ws.Range("D9")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "3.487.42"
$c.Style = $origStyle
$ws.Range("E9").Value = "  +3.14%  "
$c = $ws.Range("D10")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "7.28"
$c.Style = $origStyle
$ws.Range("E10").Value = "  -1.57%  "
$ws.Range("E11").Value = "  +3.67%  "
$ws.Range("E12").Value = "  +3.51%  "
$c = $ws.Range("D13")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "4.089.77"
$c.Style = $origStyle
$ws.Range("E13").Value = "  +3.15%  "
$ws.Range("E14").Value = "  +0.55%  "
$ws.Range("E15").Value = "  +2.78%  "
$c = $ws.Range("D16")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "28.71"
$c.Style = $origStyle
$ws.Range("E16").Value = "  +6.15%  "
$c = $ws.Range("D17")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "65.327.30"
$c.Style = $origStyle
$ws.Range("E17").Value = "  +3.65%  "
$c = $ws.Range("D18")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "3.529.64"
$c.Style = $origStyle
$ws.Range("E18").Value = "  +4.32%  "
$c = $ws.Range("D19")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "6.45"
$c.Style = $origStyle
$ws.Range("E19").Value = "  +3.38%  "
$c = $ws.Range("D20")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "14.38"
$c.Style = $origStyle
$ws.Range("E20").Value = "  +2.43%  "
$c = $ws.Range("D21")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "384.36"
$c.Style = $origStyle
$ws.Range("E21").Value = "  +2.16%  "
$c = $ws.Range("D22")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "8.24"
$c.Style = $origStyle
$ws.Range("E22").Value = "  +3.03%  "
$c = $ws.Range("D23")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.552"
$c.Style = $origStyle
$ws.Range("E23").Value = "  +4.57%  "
$c = $ws.Range("D24")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "72.84"
$c.Style = $origStyle
$ws.Range("E24").Value = "  +2.28%  "
$c = $ws.Range("D25")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = $origStyle
$ws.Range("E25").Value = "  +0.25%  "
$ws.Range("E26").Value = "  +2.44%  "
$c = $ws.Range("D27")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "10.18"
$c.Style = $origStyle
$ws.Range("E27").Value = "  +8.23%  "
$ws.Range("E28").Value = "  +1.25%  "
$c = $ws.Range("D29")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = $origStyle
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("E30").Value = "  +12.96%  "
$c = $ws.Range("D31")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "6.21"
$c.Style = $origStyle
$ws.Range("E31").Value = "  +2.61%  "
$ws.Range("E32").Value = "  +3.56%  "
$c = $ws.Range("D33")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "23.73"
$c.Style = $origStyle
$ws.Range("E33").Value = "  +2.67%  "
$c = $ws.Range("D34")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "7.19"
$c.Style = $origStyle
$ws.Range("E34").Value = "  +6.37%  "
$c = $ws.Range("D35")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.61"
$c.Style = $origStyle
$ws.Range("E35").Value = "  +12.15%  "
$c = $ws.Range("D36")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "162.58"
$c.Style = $origStyle
$ws.Range("E36").Value = "  +1.85%  "
$ws.Range("B37").Value = "Stacks"
$ws.Range("C37").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$c = $ws.Range("D37")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.93"
$c.Style = $origStyle
$ws.Range("E37").Value = "  +5.99%  "
$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$c = $ws.Range("D38")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "3.022.55"
$c.Style = $origStyle
$ws.Range("E38").Value = "  +2.10%  "
$c = $ws.Range("D39")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.0777"
$c.Style = $origStyle
$ws.Range("E39").Value = "  +2.27%  "
$ws.Range("E40").Value = "  -0.23%  "
$c = $ws.Range("D41")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "4.59"
$c.Style = $origStyle
$ws.Range("E41").Value = "  +6.61%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Range("D42")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "6.64"
$c.Style = $origStyle
$ws.Range("E42").Value = "  +4.36%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c = $ws.Range("D43")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.0320"
$c.Style = $origStyle
$ws.Range("E43").Value = "  +0.90%  "
$c = $ws.Range("D44")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "42.96"
$c.Style = $origStyle
$ws.Range("E44").Value = "  +3.19%  "
$c = $ws.Range("D45")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.781"
$c.Style = $origStyle
$ws.Range("E45").Value = "  +4.13%  "
$c = $ws.Range("D46")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "25.85"
$c.Style = $origStyle
$ws.Range("E46").Value = "  +10.96%  "
$ws.Range("E47").Value = "  +5.07%  "
$c = $ws.Range("D48")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "317.70"
$c.Style = $origStyle
$ws.Range("E48").Value = "  +10.58%  "
$c = $ws.Range("D49")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.110"
$c.Style = $origStyle
$ws.Range("E49").Value = "  +7.44%  "
$ws.Range("E50").Value = "  +6.30%  "
$ws.Range("B51").Value = "SuiNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$c = $ws.Range("D51")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.885"
$c.Style = $origStyle
$ws.Range("E51").Value = "  +6.65%  "
